$wb = $excel.ActiveWorkbook

# --- Update the conversion message on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.51 = 38803.02 pesos`n✅ 38803.02 pesos = 9.47 = 956.63 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 105.198
$wsTasas.Range("O10").Value = 4082
$wsTasas.Range("N12").Value = 4098
$wsTasas.Range("O12").Value = 101.03
